$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.255029439926147
$ws.Range("B1").Value = 2.539771795272827
$ws.Range("C1").Value = 3.762845993041992
$ws.Range("D1").Value = 2.725200891494751
$ws.Range("E1").Value = 1.072168111801147
